$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D9").Value = "[채용공고] 기사 형태의 전문 보고서 작성 업무 – 인턴레벨부터 훈련 진행"
$ws.Range("E9").Value = "https://pdsi.pabii.com/notice-hiring-special-reports-202303/#utm_source=rss&utm_medium=rss&utm_campaign=notice-hiring-special-reports-202303"

$ws.Range("D32").Value = "Bundling!!! Tensorflow model에 필요한 패키지나 모듈을 같이 저장시키는 방법!"
$ws.Range("E32").Value = "https://dodonam.tistory.com/412"

$ws.Range("D44").Value = "GPT-4 소개 및 간단 요약"
$ws.Range("E44").Value = "https://engineering-ladder.tistory.com/124"

$ws.Range("D46").Value = "갱년기증후군 (폐경기증후군)"
$ws.Range("E46").Value = "https://bioinformaticsandme.tistory.com/528"

$ws.Range("D51").Value = "[PostgreSQL] postgresql-client 설치 후 터미널에서 데이터베이스 접속 명령어"
$ws.Range("E51").Value = "https://bskyvision.com/entry/PostgreSQL-postgresql-client-%EC%84%A4%EC%B9%98-%ED%9B%84-%ED%84%B0%EB%AF%B8%EB%84%90%EC%97%90%EC%84%9C-%EB%8D%B0%EC%9D%B4%ED%84%B0%EB%B2%A0%EC%9D%B4%EC%8A%A4-%EC%A0%91%EC%86%8D-%EB%AA%85%EB%A0%B9%EC%96%B4"
